$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2026-02-26", "substance active", 44, 1),
    @("2026-02-26", "substance active", 45, 3),
    @("2026-02-26", "eaux de surface", 104, 1),
    @("2026-02-26", "zone tampon", 105, 2),
    @("2026-02-26", "eaux de surface", 106, 1),
    @("2026-02-26", "zone tampon", 108, 1),
    @("2026-02-26", "développement durable", 171, 1),
    @("2026-02-26", "ruissellement", 218, 1)
)

$startRow = 481
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
